# Auto-generated script applying cell-level updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 246.8
$ws.Range("I5").Value = 58.5
$ws.Range("K5").Value = 58.5
$ws.Range("M5").Value = 56.5
$ws.Range("H40").Value = 140162.64
$ws.Range("I40").Value = 752919.5
$ws.Range("J40").Value = 3994.4443
$ws.Range("K40").Value = 752919.5
$ws.Range("L40").Value = 3994.4443
$ws.Range("M40").Value = -752744.5
$ws.Range("N40").Value = -4344.4443
$ws.Range("H86").Value = 132359064
$ws.Range("I86").Value = 111118950
$ws.Range("J86").Value = 156254190
$ws.Range("K86").Value = 111118950
$ws.Range("L86").Value = 156254190
$ws.Range("M86").Value = -111117827
$ws.Range("N86").Value = -156256436
$ws.Range("H87").Value = 164988
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 164988
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 164988
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -167484
$ws.Range("H89").Value = 132359064
$ws.Range("I89").Value = 111118950
$ws.Range("J89").Value = 156254190
$ws.Range("K89").Value = 555594750
$ws.Range("L89").Value = 781270950
$ws.Range("M89").Value = -555589134
$ws.Range("N89").Value = -781282182
$ws.Range("H90").Value = 164988
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 164988
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 494964
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -507444
$ws.Range("H92").Value = 33334128
$ws.Range("I92").Value = 41667548
$ws.Range("J92").Value = 443.33334
$ws.Range("K92").Value = 41667548
$ws.Range("L92").Value = 443.33334
$ws.Range("M92").Value = -41666300
$ws.Range("N92").Value = -2939.33334
$ws.Range("H132").Value = 4912.1816
$ws.Range("I132").Value = 5054.5366
$ws.Range("K132").Value = 15163.6098
$ws.Range("M132").Value = -12633.6098
$ws.Range("H138").Value = 3755.87
$ws.Range("I138").Value = 1872.7273
$ws.Range("J138").Value = 3988.618
$ws.Range("K138").Value = 5618.1819
$ws.Range("L138").Value = 11965.854
$ws.Range("M138").Value = -478.1818999999996
$ws.Range("N138").Value = -22245.854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 447
$ws.Range("I5").Value = 471.25
$ws.Range("K5").Value = 471.25
$ws.Range("M5").Value = -359.25
$ws.Range("H61").Value = 2470.8154
$ws.Range("I61").Value = 2263.2068
$ws.Range("J61").Value = 4191
$ws.Range("K61").Value = 2263.2068
$ws.Range("L61").Value = 4191
$ws.Range("M61").Value = -2051.2068
$ws.Range("N61").Value = -4615
$ws.Range("H102").Value = 1859.3572
$ws.Range("I102").Value = 1575.6364
$ws.Range("J102").Value = 2899.6667
$ws.Range("K102").Value = 1575.6364
$ws.Range("L102").Value = 2899.6667
$ws.Range("M102").Value = 46.36359999999991
$ws.Range("N102").Value = -6143.6667
$ws.Range("H132").Value = 113211.22
$ws.Range("I132").Value = 136837.19
$ws.Range("J132").Value = 3941.125
$ws.Range("K132").Value = 410511.57
$ws.Range("L132").Value = 11823.375
$ws.Range("M132").Value = -407981.57
$ws.Range("N132").Value = -16883.375
$ws.Range("H136").Value = 2470.8154
$ws.Range("I136").Value = 2263.2068
$ws.Range("J136").Value = 4191
$ws.Range("K136").Value = 6789.6204
$ws.Range("L136").Value = 12573
$ws.Range("M136").Value = -4239.6204
$ws.Range("N136").Value = -17673

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 447
$ws.Range("I4").Value = 471.25
$ws.Range("K4").Value = 471.25
$ws.Range("M4").Value = -356.25
$ws.Range("H86").Value = 2043.8928
$ws.Range("I86").Value = 1791.1333
$ws.Range("J86").Value = 2335.5386
$ws.Range("K86").Value = 1791.1333
$ws.Range("L86").Value = 2335.5386
$ws.Range("M86").Value = -668.1333
$ws.Range("N86").Value = -4581.5386
$ws.Range("H89").Value = 2043.8928
$ws.Range("I89").Value = 1791.1333
$ws.Range("J89").Value = 2335.5386
$ws.Range("K89").Value = 8955.666499999999
$ws.Range("L89").Value = 11677.693
$ws.Range("M89").Value = -3339.666499999999
$ws.Range("N89").Value = -22909.693
$ws.Range("H131").Value = 41784.5
$ws.Range("J131").Value = 41784.5
$ws.Range("L131").Value = 41784.5
$ws.Range("N131").Value = -51864.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6550000
$ws.Range("I6").Value = 6550000
$ws.Range("K6").Value = 6550000
$ws.Range("M6").Value = -6549887
$ws.Range("H7").Value = 169.375
$ws.Range("I7").Value = 45.545456
$ws.Range("J7").Value = 441.8
$ws.Range("K7").Value = 45.545456
$ws.Range("L7").Value = 441.8
$ws.Range("M7").Value = 67.454544
$ws.Range("N7").Value = -667.8
$ws.Range("H22").Value = 614.4666999999999
$ws.Range("I22").Value = 570
$ws.Range("J22").Value = 681.1667
$ws.Range("K22").Value = 570
$ws.Range("L22").Value = 681.1667
$ws.Range("M22").Value = -220
$ws.Range("N22").Value = -1381.1667
$ws.Range("H25").Value = 4250
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -5348
$ws.Range("H31").Value = 2386.4888
$ws.Range("I31").Value = 1734.3334
$ws.Range("J31").Value = 3364.7222
$ws.Range("K31").Value = 1734.3334
$ws.Range("L31").Value = 3364.7222
$ws.Range("M31").Value = -1439.3334
$ws.Range("N31").Value = -3954.7222
$ws.Range("H34").Value = 2386.4888
$ws.Range("I34").Value = 1734.3334
$ws.Range("J34").Value = 3364.7222
$ws.Range("K34").Value = 1734.3334
$ws.Range("L34").Value = 3364.7222
$ws.Range("M34").Value = -1532.3334
$ws.Range("N34").Value = -3768.7222
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H62").Value = 3232.3333
$ws.Range("I62").Value = 3232.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3232.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2608.3333
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3232.3333
$ws.Range("I65").Value = 3232.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16161.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13041.6665
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 30285
$ws.Range("I74").Value = 30285
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 30285
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -29411
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 30285
$ws.Range("I77").Value = 30285
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 90855
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -86487
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1257.381
$ws.Range("I131").Value = 923.46155
$ws.Range("J131").Value = 1800
$ws.Range("K131").Value = 2770.38465
$ws.Range("L131").Value = 5400
$ws.Range("M131").Value = 2269.61535
$ws.Range("N131").Value = -15480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 32748.762
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 32748.762
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 98246.28599999999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -103346.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2966.6
$ws.Range("I46").Value = 1060.2
$ws.Range("J46").Value = 3919.8
$ws.Range("K46").Value = 1060.2
$ws.Range("L46").Value = 3919.8
$ws.Range("M46").Value = -872.2
$ws.Range("N46").Value = -4295.8
$ws.Range("H68").Value = 13970.3
$ws.Range("I68").Value = 29666.666
$ws.Range("J68").Value = 7243.2856
$ws.Range("K68").Value = 29666.666
$ws.Range("L68").Value = 7243.2856
$ws.Range("M68").Value = -28917.666
$ws.Range("N68").Value = -8741.285599999999
$ws.Range("H71").Value = 13970.3
$ws.Range("I71").Value = 29666.666
$ws.Range("J71").Value = 7243.2856
$ws.Range("K71").Value = 148333.33
$ws.Range("L71").Value = 36216.428
$ws.Range("M71").Value = -144589.33
$ws.Range("N71").Value = -43704.428
$ws.Range("H132").Value = 1113365.5
$ws.Range("I132").Value = 1668381.9
$ws.Range("K132").Value = 5005145.699999999
$ws.Range("M132").Value = -5002615.699999999
$ws.Range("H136").Value = 2326.074
$ws.Range("I136").Value = 1495.25
$ws.Range("K136").Value = 4485.75
$ws.Range("M136").Value = -1935.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 86665.664
$ws.Range("J130").Value = 86665.664
$ws.Range("L130").Value = 86665.664
$ws.Range("N130").Value = -96705.664
$ws.Range("H136").Value = 50858.855
$ws.Range("I136").Value = 2707.9092
$ws.Range("K136").Value = 8123.7276
$ws.Range("M136").Value = -5573.7276
